# Case_3_235 (380 kV) results: update vm_pu.xlsx bus voltage magnitudes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{}
$data[2] = @(1.02, 1.034549298079764, 1.039237236565296, 1.033693968881767, 1.044472820629332, 1.040201713768006, 1.039667933892078, 1.042022995502928, 1.036495582016846, 1.047243773065315, 1.017095310108538)
$data[3] = @(1.02, 1.035592239019622, 1.040073052189016, 1.034584819740886, 1.046789062068997, 1.0406002617218, 1.040353580208111, 1.042668920881376, 1.037195247949527, 1.049367338425143, 1.017331119084254)
$data[4] = @(1.02, 1.036266555290422, 1.04061324969089, 1.035161092158067, 1.048281843765013, 1.040856171610722, 1.040796136034545, 1.04308560366257, 1.037647180220787, 1.050735140636954, 1.017483138216689)
$data[5] = @(1.02, 1.036549911944327, 1.040840199626356, 1.035403318487171, 1.048908013170412, 1.040963285381175, 1.040981924385372, 1.043260474210066, 1.03783698307613, 1.051308693226591, 1.017546912365824)
$data[6] = @(1.02, 1.036597481461356, 1.040878296833261, 1.03544398713463, 1.049013068940392, 1.040981242725054, 1.041013103753282, 1.043289818043941, 1.037868840747661, 1.051404909909178, 1.017557612444137)
$data[7] = @(1.02, 1.036270342006922, 1.040616282789894, 1.035164328947389, 1.048290216108415, 1.040857604716192, 1.040798619574591, 1.043087941478427, 1.037649717117821, 1.050742810211498, 1.017483990898366)
$data[8] = @(1.02, 1.034901877389886, 1.039519836252706, 1.033995071754384, 1.045256865533383, 1.040336815914951, 1.039899880903944, 1.042241554123393, 1.036732203712321, 1.047962763933348, 1.017175120203566)
$data[9] = @(1.02, 1.032486247928765, 1.03758283242321, 1.031933336031656, 1.039864275285698, 1.039403856515247, 1.03830764852633, 1.040740255428783, 1.035109237109187, 1.043014298616173, 1.016626496434572)
$data[10] = @(1.02, 1.030872828646408, 1.036288052610002, 1.030557830541328, 1.036235022550647, 1.0387714513984, 1.037240287407188, 1.039732622321974, 1.034022986701529, 1.039679781627275, 1.016257781104034)
$data[11] = @(1.02, 1.030173451286913, 1.035726550840306, 1.029961956254317, 1.034654888983438, 1.038495099728783, 1.036776684930313, 1.039294668413303, 1.033551589128294, 1.038226986262501, 1.016097410695983)
$data[12] = @(1.02, 1.029913554423105, 1.035517853042596, 1.029740579013128, 1.034066613582789, 1.038392068947211, 1.036604265043024, 1.039131743197967, 1.033376331903072, 1.037685970158404, 1.016037733850601)
$data[13] = @(1.02, 1.0299693085306, 1.035562625404113, 1.029788067165117, 1.03419286217584, 1.038414186717938, 1.036641259571255, 1.039166702553746, 1.03341393241389, 1.037802083224465, 1.016050539643378)
$data[14] = @(1.02, 1.03015197053698, 1.035709302504847, 1.029943658032802, 1.034606289623685, 1.038486590974267, 1.03676243709722, 1.039281206073761, 1.033537105572228, 1.038182294165315, 1.01609248000012)
$data[15] = @(1.02, 1.030264499025925, 1.03579965765865, 1.030039516935883, 1.034860836543471, 1.038531150973694, 1.036837069695141, 1.039351722304146, 1.033612975457126, 1.038416370116084, 1.016118306484004)
$data[16] = @(1.02, 1.030919227800826, 1.036325299411926, 1.030597370849055, 1.036339704895466, 1.038789738644298, 1.03727102485685, 1.039761653034729, 1.034054249612651, 1.039776007213974, 1.016268409241335)
$data[17] = @(1.02, 1.031329717057899, 1.036654790114933, 1.030947223665246, 1.037265014470922, 1.038951267813848, 1.037542849217016, 1.040018350063551, 1.034330768098694, 1.040626452108633, 1.016362372934071)
$data[18] = @(1.02, 1.031569075769508, 1.036846894304215, 1.031151260869003, 1.037803901021419, 1.039045242517432, 1.037701262167863, 1.040167918857706, 1.034491956085434, 1.041121643723984, 1.01641711151013)
$data[19] = @(1.02, 1.031650678669367, 1.036912383015427, 1.031220827934838, 1.037987507826093, 1.039077244429196, 1.037755253632651, 1.040218891198327, 1.034546900034851, 1.041290346909754, 1.016435764295478)
$data[20] = @(1.02, 1.031285682992379, 1.036619447399564, 1.030909690477135, 1.037165823773878, 1.038933962373378, 1.03751369928804, 1.039990825290326, 1.034301110701179, 1.040535296520503, 1.01635229864723)
$data[21] = @(1.02, 1.030098184385194, 1.035666113394658, 1.029897841605483, 1.034484582994712, 1.038465280283828, 1.036726759356471, 1.039247494555932, 1.033500838560161, 1.038070370044099, 1.016080132602212)
$data[22] = @(1.02, 1.029350878399133, 1.035065954878833, 1.029261402130805, 1.032790984951957, 1.038168392443023, 1.036230720445967, 1.038778687224091, 1.032996753558495, 1.036512548716105, 1.01590838494528)
$data[23] = @(1.02, 1.029747104893305, 1.035384183215227, 1.029598815188066, 1.033689547446645, 1.038325988844574, 1.03649380028006, 1.039027348897077, 1.033264066755193, 1.037339154087693, 1.015999491252853)
$data[24] = @(1.02, 1.031305580316252, 1.036635417505505, 1.03092665019001, 1.037210646345773, 1.038941782706472, 1.037526871311609, 1.040003263039713, 1.034314511913536, 1.040576488458058, 1.01635685099667)
$data[25] = @(1.02, 1.03311126266933, 1.038084191433756, 1.032466515259198, 1.041264245037345, 1.039646874949653, 1.038720304117472, 1.041129559413813, 1.035529557721998, 1.044299700193808, 1.016768848741715)

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $vals[$i]
    }
}

